$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.750.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3
$ws.Range("D3").Value = "'1.868.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "

# Row 4
$ws.Range("D4").Value = "'0.9982"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'0.7296"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "

# Row 6
$ws.Range("D6").Value = "'240.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7
$ws.Range("D7").Value = "'0.9991"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.95%  "

# Row 9
$ws.Range("D9").Value = "'0.07096"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.42%  "

# Row 10
$ws.Range("D10").Value = "'24.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.58%  "

# Row 11
$ws.Range("D11").Value = "'0.08193"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.56%  "

# Row 12
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7379"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.875.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "

# Row 14
$ws.Range("D14").Value = "'5.323"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'92.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16
$ws.Range("D16").Value = "'29.778.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "

# Row 17
$ws.Range("D17").Value = "'6.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18
$ws.Range("D18").Value = "'247.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.81%  "

# Row 19
$ws.Range("E19").Value = "  -1.74%  "

# Row 20
$ws.Range("D20").Value = "'0.000007784"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.164.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "

# Row 23
$ws.Range("D23").Value = "'0.9983"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").Value = "'7.750"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.62%  "

# Row 25
$ws.Range("D25").Value = "'0.1543"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.54%  "

# Row 26
$ws.Range("D26").Value = "'9.167"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.48%  "

# Row 27
$ws.Range("D27").Value = "'163.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("D28").Value = "'18.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.64%  "

# Row 29
$ws.Range("D29").Value = "'2.005"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.98%  "

# Row 30
$ws.Range("D30").Value = "'1.442"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.01%  "

# Row 31
$ws.Range("E31").Value = "  -2.91%  "

# Row 32
$ws.Range("D32").Value = "'1.518"
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "  -3.16%  "

# Row 34
$ws.Range("D34").Value = "'0.05279"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "

# Row 35
$ws.Range("E35").Value = "  -0.40%  "

# Row 36
$ws.Range("D36").Value = "'0.7426"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "

# Row 37
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("D38").Value = "'2.687"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "

# Row 39
$ws.Range("D39").Value = "'0.01927"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "

# Row 40
$ws.Range("D40").Value = "'2.728"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.60%  "

# Row 41
$ws.Range("D41").Value = "'0.4442"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.64%  "

# Row 42
$ws.Range("D42").Value = "'5.982"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.71%  "

# Row 43
$ws.Range("D43").Value = "'0.8658"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "

# Row 44
$ws.Range("D44").Value = "'71.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.57%  "

# Row 45
$ws.Range("D45").Value = "'1.038.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.43%  "

# Row 46
$ws.Range("D46").Value = "'0.9985"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
$ws.Range("D47").Value = "'103.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "

# Row 48
$ws.Range("D48").Value = "'7.428"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.40%  "

# Row 49
$ws.Range("D49").Value = "'1.809"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.92%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "'2.054.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.08%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.499"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
